$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.997.41'
$ws.Range("E2").Value = '  +1.58%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.252.65'
$ws.Range("E3").Value = '  +0.93%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '271.65'
$ws.Range("E5").Value = '  +5.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.31'
$ws.Range("E6").Value = '  +11.31%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  +1.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.31'
$ws.Range("E10").Value = '  +5.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0926'
$ws.Range("E11").Value = '  +0.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.59'
$ws.Range("E12").Value = '  +7.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.105'
$ws.Range("E13").Value = '  +2.16%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.590.37'
$ws.Range("E14").Value = '  +1.06%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.90'
$ws.Range("E15").Value = '  +1.63%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.249.28'
$ws.Range("E16").Value = '  +0.51%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.797'
$ws.Range("E17").Value = '  +0.44%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.944.98'
$ws.Range("E18").Value = '  +1.66%  '

$ws.Range("E19").Value = '  -0.51%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.01'
$ws.Range("E20").Value = '  -0.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.51'
$ws.Range("E21").Value = '  -1.15%  '

$ws.Range("E22").Value = '  +2.98%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.90'
$ws.Range("E23").Value = '  +1.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.90'
$ws.Range("E24").Value = '  -4.29%  '

$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.51'
$ws.Range("E26").Value = '  +13.18%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.85'
$ws.Range("E27").Value = '  +0.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.54'
$ws.Range("E28").Value = '  +6.28%  '

$ws.Range("E29").Value = '  +5.35%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '39.78'
$ws.Range("E30").Value = '  -5.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.27'
$ws.Range("E31").Value = '  +0.95%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.90'
$ws.Range("E32").Value = '  +2.11%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0896'
$ws.Range("E33").Value = '  +2.97%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.38'
$ws.Range("E34").Value = '  +2.22%  '

$ws.Range("E35").Value = '  +0.93%  '

$ws.Range("E36").Value = '  +3.42%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0353'
$ws.Range("E37").Value = '  -4.53%  '

$ws.Range("E38").Value = '  -1.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.47'
$ws.Range("E39").Value = '  +19.01%  '

$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.64'
$ws.Range("E40").Value = '  -4.35%  '

$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.20'
$ws.Range("E41").Value = '  +2.64%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '64.68'
$ws.Range("E42").Value = '  +4.66%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.43'
$ws.Range("E43").Value = '  +1.62%  '

$ws.Range("E44").Value = '  -0.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.53'
$ws.Range("E45").Value = '  -0.92%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0988'
$ws.Range("E46").Value = '  +0.73%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.99'
$ws.Range("E47").Value = '  -2.66%  '

$ws.Range("E48").Value = '  +4.73%  '

$ws.Range("E49").Value = '  +1.41%  '

$ws.Range("E50").Value = '  +2.18%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.431'
$ws.Range("E51").Value = '  -9.52%  '
